$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting -
# the source cells are inline strings like '42.575.84' / '1.00' that
# would otherwise be auto-coerced to numbers (losing trailing zeros
# and the multi-dot 'thousands' formatting) by the normal Value setter.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '42.718.76'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '2.311.87'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '303.11'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').Value = '100.12'
$ws.Range('E6').Value = '  -3.81%  '
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  -3.44%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.504'
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('D10').Value = '34.87'
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '6.75'
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = '2.666.76'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '15.76'
$ws.Range('E15').Value = '  +4.52%  '
$ws.Range('D16').Value = '2.298.39'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '0.807'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '42.595.03'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('D19').Value = '0.0₃0908'
$ws.Range('E19').Value = '  -1.44%  '
$ws.Range('D20').Value = '6.09'
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').Value = '11.53'
$ws.Range('E21').Value = '  -3.54%  '
$ws.Range('D22').Value = '67.96'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').Value = '235.73'
$ws.Range('E23').Value = '  -1.80%  '
$ws.Range('D24').Value = '1.97'
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').Value = '  -2.82%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '25.01'
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('D28').Value = '2.18'
$ws.Range('E28').Value = '  +3.08%  '
$ws.Range('D29').Value = '34.83'
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('D30').Value = '164.62'
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('D31').Value = '9.16'
$ws.Range('E31').Value = '  -4.16%  '
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').Value = '5.03'
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('E34').Value = '  -4.84%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '4.46'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').Value = '16.81'
$ws.Range('E36').Value = '  -7.60%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.0702'
$ws.Range('E37').Value = '  -4.44%  '
$ws.Range('D38').Value = '2.90'
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('D39').Value = '1.81'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('E40').Value = '  -5.85%  '
$ws.Range('E41').Value = '  -3.00%  '
$ws.Range('D42').Value = '2.46'
$ws.Range('E42').Value = '  -4.93%  '
$ws.Range('D43').Value = '1.970.68'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').Value = '0.0281'
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').Value = '18.60'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '10.22'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').Value = '2.90'
$ws.Range('E47').Value = '  -6.02%  '
$ws.Range('D48').Value = '55.70'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').Value = '2.87'
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('D50').Value = '2.532.99'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '4.69'
$ws.Range('E51').Value = '  +0.61%  '
